# Atualização de bases das ligas, do dia: 24-02-2024 às 12:40
#
# Rows 98 and 99 swap their full data (all columns except the running id in
# column A, which stays tied to the row position, i.e. to the order in
# which rows appear on the sheet).
# Rows 102, 103 and 104 rotate (102<-104, 103<-102, 104<-103), again
# keeping column A fixed per row.
# Rows 121 and 122 swap their full data, same rule.
# The last data row (210) is removed entirely, shrinking the used range
# from A1:AC210 down to A1:AC209.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-MatchRow($Row, $Id, $Home, $Away, $Fthg, $Ftag, $Ftr, $K, $L, $M, $N, $O, $P, $Q, $R, $S, $T, $U, $V, $W, $X, $Y, $Z, $AA, $AB, $AC) {
    $ws.Cells.Item($Row, 2).Value = $Id      # B  id
    $ws.Cells.Item($Row, 6).Value = $Home    # F  HomeTeam
    $ws.Cells.Item($Row, 7).Value = $Away    # G  AwayTeam
    $ws.Cells.Item($Row, 8).Value = $Fthg    # H  FTHG
    $ws.Cells.Item($Row, 9).Value = $Ftag    # I  FTAG
    $ws.Cells.Item($Row, 10).Value = $Ftr    # J  FTR
    $ws.Cells.Item($Row, 11).Value = $K      # K  oddH_op
    $ws.Cells.Item($Row, 12).Value = $L      # L  oddD_op
    $ws.Cells.Item($Row, 13).Value = $M      # M  oddA_op
    $ws.Cells.Item($Row, 14).Value = $N      # N  oddH
    $ws.Cells.Item($Row, 15).Value = $O      # O  oddD
    $ws.Cells.Item($Row, 16).Value = $P      # P  oddA
    $ws.Cells.Item($Row, 17).Value = $Q      # Q  Ah
    $ws.Cells.Item($Row, 18).Value = $R      # R  oddAHH
    $ws.Cells.Item($Row, 19).Value = $S      # S  oddAHA
    $ws.Cells.Item($Row, 20).Value = $T      # T  AhOU
    $ws.Cells.Item($Row, 21).Value = $U      # U  oddAHOver
    $ws.Cells.Item($Row, 22).Value = $V      # V  oddAHUnder
    $ws.Cells.Item($Row, 23).Value = $W      # W  PLH
    $ws.Cells.Item($Row, 24).Value = $X      # X  PLD
    $ws.Cells.Item($Row, 25).Value = $Y      # Y  PLA
    $ws.Cells.Item($Row, 26).Value = $Z      # Z  PL_Ahh
    $ws.Cells.Item($Row, 27).Value = $AA     # AA PL_Aha
    $ws.Cells.Item($Row, 28).Value = $AB     # AB PL_AhOver
    $ws.Cells.Item($Row, 29).Value = $AC     # AC PL_AhUnder
}

# ---- Row 98 <- former row 99 data ----
Set-MatchRow 98 "6425845" "NK Kustosija" "NK Dubrava Zagreb" "0" "1" "A" 2.1 3.5 2.875 2.375 3.5 2.5 0 1.825 1.975 2.5 2 1.8 -1 -1 1.5 -1 0.9750000000000001 -1 0.8

# ---- Row 99 <- former row 98 data ----
Set-MatchRow 99 "6425846" "Orijent" "NK Hrvatski Dragovoljac" "2" "2" "D" 1.444 4 6 1.25 5.25 9 -1.75 1.95 1.85 3 1.8 2 -1 4.25 -1 -1 0.8500000000000001 0.8 -1

# ---- Row 102 <- former row 104 data ----
Set-MatchRow 102 "6834733" "HNK Cibalia" "NK Croatia Zmijavci" "1" "0" "H" 1.65 3.5 4.5 1.909 3.3 3.3 -0.5 2 1.8 2.25 1.95 1.85 0.909 -1 -1 1 -1 -1 0.8500000000000001

# ---- Row 103 <- former row 102 data ----
Set-MatchRow 103 "6834732" "Orijent" "NK Dubrava Zagreb" "1" "1" "D" 2.1 3.2 3.1 2.3 3.1 2.8 -0.25 2.025 1.775 2.25 1.825 1.975 -1 2.1 -1 -0.5 0.3875 -0.5 0.4875

# ---- Row 104 <- former row 103 data ----
Set-MatchRow 104 "6834729" "NK Solin" "Bijelo Brdo" "0" "0" "D" 2.1 3.2 3.1 2.05 3.25 3.25 -0.25 1.8 2 2.5 2 1.8 -1 2.25 -1 -0.5 0.5 -1 0.8

# ---- Row 121 <- former row 122 data ----
Set-MatchRow 121 "6834749" "Vukovar 91" "HNK Cibalia" "3" "1" "H" 2.05 3.1 3.25 1.95 3.2 3.5 -0.5 2 1.8 2.25 1.925 1.875 0.95 -1 -1 1 -1 0.925 -1

# ---- Row 122 <- former row 121 data ----
Set-MatchRow 122 "6834750" "NK Solin" "Orijent" "4" "1" "H" 2.1 3.3 3 2.1 3.3 3 -0.25 1.875 1.925 2.5 1.825 1.975 1.1 -1 -1 0.875 -1 0.825 -1

# ---- Remove the last data row entirely (used range shrinks to AC209) ----
$ws.Rows(210).Delete()
